$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value2 = 46616910  # was 46616996
$ws.Range("J64").Value2 = 45458330  # was 45458428
$ws.Range("L64").Value2 = 45458330  # was 45458428
$ws.Range("N64").Value2 = -45458826  # was -45458924
$ws.Range("H67").Value2 = 46616910  # was 46616996
$ws.Range("J67").Value2 = 45458330  # was 45458428
$ws.Range("L67").Value2 = 45458330  # was 45458428
$ws.Range("N67").Value2 = -45460046  # was -45460144
$ws.Range("H76").Value2 = 4655.857  # was 4623.5
$ws.Range("I76").Value2 = 4655.857  # was 4623.5
$ws.Range("K76").Value2 = 4655.857  # was 4623.5
$ws.Range("M76").Value2 = -4340.857  # was -4308.5
$ws.Range("H79").Value2 = 4655.857  # was 4623.5
$ws.Range("I79").Value2 = 4655.857  # was 4623.5
$ws.Range("K79").Value2 = 4655.857  # was 4623.5
$ws.Range("M79").Value2 = -3563.857  # was -3531.5
$ws.Range("H101").Value2 = 239  # was 2446.1
$ws.Range("I101").Value2 = 245.2  # was 1179.7142
$ws.Range("J101").Value2 = 208  # was 5401
$ws.Range("K101").Value2 = 735.5999999999999  # was 3539.1426
$ws.Range("L101").Value2 = 624  # was 16203
$ws.Range("M101").Value2 = 886.4000000000001  # was -1917.1426
$ws.Range("N101").Value2 = -3868  # was -19447
$ws.Range("H132").Value2 = 13481.131  # was 10619.848
$ws.Range("I132").Value2 = 6327.5557  # was 2986.524
$ws.Range("J132").Value2 = 15221.189  # was 14838.263
$ws.Range("K132").Value2 = 18982.6671  # was 8959.572
$ws.Range("L132").Value2 = 45663.567  # was 44514.789
$ws.Range("M132").Value2 = -16452.6671  # was -6429.572
$ws.Range("N132").Value2 = -50723.567  # was -49574.789
$ws.Range("H135").Value2 = 1894.325  # was 2137.5557
$ws.Range("I135").Value2 = 525.82355  # was 601.2857
$ws.Range("J135").Value2 = 2905.8262  # was 3115.182
$ws.Range("K135").Value2 = 4732.41195  # was 5411.571300000001
$ws.Range("L135").Value2 = 26152.4358  # was 28036.638
$ws.Range("M135").Value2 = -2197.41195  # was -2876.571300000001
$ws.Range("N135").Value2 = -31222.4358  # was -33106.638
$ws.Range("H138").Value2 = 5622.2856  # was 5307.352
$ws.Range("I138").Value2 = 2287.5  # was 2475.2222
$ws.Range("J138").Value2 = 6477.359  # was 5873.778
$ws.Range("K138").Value2 = 6862.5  # was 7425.6666
$ws.Range("L138").Value2 = 19432.077  # was 17621.334
$ws.Range("M138").Value2 = -1722.5  # was -2285.6666
$ws.Range("N138").Value2 = -29712.077  # was -27901.334

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value2 = 16999  # was 30000000
$ws.Range("J92").Value2 = 16999  # was 30000000
$ws.Range("L92").Value2 = 16999  # was 30000000
$ws.Range("N92").Value2 = -21991  # was -30004992
$ws.Range("H97").Value2 = 1018.6  # was 1144.4286
$ws.Range("I97").Value2 = 1026.5714  # was 1168.5
$ws.Range("K97").Value2 = 1026.5714  # was 1168.5
$ws.Range("M97").Value2 = -530.5714  # was -672.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value2 = 341.2  # was 326.9
$ws.Range("I80").Value2 = 531.75  # was 426.6
$ws.Range("J80").Value2 = 214.16667  # was 227.2
$ws.Range("K80").Value2 = 531.75  # was 426.6
$ws.Range("L80").Value2 = 214.16667  # was 227.2
$ws.Range("M80").Value2 = 466.25  # was 571.4
$ws.Range("N80").Value2 = -2210.16667  # was -2223.2
$ws.Range("H83").Value2 = 341.2  # was 326.9
$ws.Range("I83").Value2 = 531.75  # was 426.6
$ws.Range("J83").Value2 = 214.16667  # was 227.2
$ws.Range("K83").Value2 = 2658.75  # was 2133
$ws.Range("L83").Value2 = 1070.83335  # was 1136
$ws.Range("M83").Value2 = 2333.25  # was 2859
$ws.Range("N83").Value2 = -11054.83335  # was -11120
$ws.Range("H134").Value2 = 4299.875  # was 4494.4
$ws.Range("I134").Value2 = 3447.2942  # was 3529.3333
$ws.Range("J134").Value2 = 6370.4287  # was 7389.6
$ws.Range("K134").Value2 = 10341.8826  # was 10587.9999
$ws.Range("L134").Value2 = 19111.2861  # was 22168.8
$ws.Range("M134").Value2 = -7806.882599999999  # was -8052.999899999999
$ws.Range("N134").Value2 = -24181.2861  # was -27238.8

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value2 = 386  # was 351.5
$ws.Range("I22").Value2 = 370  # was 331.8
$ws.Range("K22").Value2 = 370  # was 331.8
$ws.Range("M22").Value2 = -20  # was 18.19999999999999
$ws.Range("H99").Value2 = 9926.434999999999  # was 9603.791999999999
$ws.Range("I99").Value2 = 12272.429  # was 11599.8
$ws.Range("K99").Value2 = 12272.429  # was 11599.8
$ws.Range("M99").Value2 = -10774.429  # was -10101.8
$ws.Range("H107").Value2 = 4546079  # was 6061189.5
$ws.Range("I107").Value2 = 9091284  # was 18181818
$ws.Range("K107").Value2 = 9091284  # was 18181818
$ws.Range("M107").Value2 = -9089364  # was -18179898
$ws.Range("H122").Value2 = 1994.9546  # was 1995.6818
$ws.Range("I122").Value2 = 1758.8125  # was 1759.8125
$ws.Range("K122").Value2 = 5276.4375  # was 5279.4375
$ws.Range("M122").Value2 = -2826.4375  # was -2829.4375
$ws.Range("H126").Value2 = 9926.434999999999  # was 9603.791999999999
$ws.Range("I126").Value2 = 12272.429  # was 11599.8
$ws.Range("K126").Value2 = 36817.287  # was 34799.39999999999
$ws.Range("M126").Value2 = -34347.287  # was -32329.39999999999
$ws.Range("H134").Value2 = 2459.6135  # was 2490.9768
$ws.Range("I134").Value2 = 2546.5134  # was 2586.389
$ws.Range("K134").Value2 = 7639.540199999999  # was 7759.167
$ws.Range("M134").Value2 = -5104.540199999999  # was -5224.167

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value2 = 6889393.5  # was 6889397.5
$ws.Range("J137").Value2 = 7872733  # was 7872738
$ws.Range("L137").Value2 = 23618199  # was 23618214
$ws.Range("N137").Value2 = -23628399  # was -23628414

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value2 = 2760.1538  # was 2793.5264
$ws.Range("I126").Value2 = 2514.423  # was 2555.32
$ws.Range("K126").Value2 = 7543.268999999999  # was 7665.960000000001
$ws.Range("M126").Value2 = -5073.268999999999  # was -5195.960000000001
$ws.Range("H132").Value2 = 4281.5903  # was 4468.035
$ws.Range("I132").Value2 = 3898.6326  # was 4055.087
$ws.Range("J132").Value2 = 5845.3335  # was 6194.909
$ws.Range("K132").Value2 = 11695.8978  # was 12165.261
$ws.Range("L132").Value2 = 17536.0005  # was 18584.727
$ws.Range("M132").Value2 = -9165.897799999999  # was -9635.261
$ws.Range("N132").Value2 = -22596.0005  # was -23644.727
$ws.Range("H140").Value2 = 71715  # was 71239.28999999999
$ws.Range("J140").Value2 = 71715  # was 71239.28999999999
$ws.Range("L140").Value2 = 71715  # was 71239.28999999999
$ws.Range("N140").Value2 = -82075  # was -81599.28999999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value2 = 2069.375  # was 3394.625
$ws.Range("I22").Value2 = 1499.2858  # was 2747.5
$ws.Range("J22").Value2 = 2512.7778  # was 3610.3333
$ws.Range("K22").Value2 = 1499.2858  # was 2747.5
$ws.Range("L22").Value2 = 2512.7778  # was 3610.3333
$ws.Range("M22").Value2 = -1204.2858  # was -2452.5
$ws.Range("N22").Value2 = -3102.7778  # was -4200.3333
$ws.Range("H27").Value2 = 2069.375  # was 3394.625
$ws.Range("I27").Value2 = 1499.2858  # was 2747.5
$ws.Range("J27").Value2 = 2512.7778  # was 3610.3333
$ws.Range("K27").Value2 = 1499.2858  # was 2747.5
$ws.Range("L27").Value2 = 2512.7778  # was 3610.3333
$ws.Range("M27").Value2 = -1392.2858  # was -2640.5
$ws.Range("N27").Value2 = -2726.7778  # was -3824.3333
$ws.Range("H82").Value2 = 3126649.8  # was 2842454.2
$ws.Range("I82").Value2 = 4466000  # was 3907812.5
$ws.Range("K82").Value2 = 4466000  # was 3907812.5
$ws.Range("M82").Value2 = -4465639  # was -3907451.5
$ws.Range("H85").Value2 = 3126649.8  # was 2842454.2
$ws.Range("I85").Value2 = 4466000  # was 3907812.5
$ws.Range("K85").Value2 = 4466000  # was 3907812.5
$ws.Range("M85").Value2 = -4464752  # was -3906564.5
$ws.Range("H132").Value2 = 3848.1128  # was 3952.4746
$ws.Range("I132").Value2 = 2789.0408  # was 2853.8262
$ws.Range("K132").Value2 = 8367.1224  # was 8561.4786
$ws.Range("M132").Value2 = -5837.1224  # was -6031.4786
$ws.Range("H136").Value2 = 2594.1924  # was 2731.6667
$ws.Range("I136").Value2 = 1974  # was 2025.75
$ws.Range("J136").Value2 = 5199  # was 6261.25
$ws.Range("K136").Value2 = 5922  # was 6077.25
$ws.Range("L136").Value2 = 15597  # was 18783.75
$ws.Range("M136").Value2 = -3372  # was -3527.25
$ws.Range("N136").Value2 = -20697  # was -23883.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value2 = 20000  # was 0
$ws.Range("I58").Value2 = 20000  # was 0
$ws.Range("K58").Value2 = 20000  # was 0
$ws.Range("M58").Value2 = -19692  # new cell
$ws.Range("H75").Value2 = 120000  # was 50000
$ws.Range("I75").Value2 = 0  # was 50000
$ws.Range("J75").Value2 = 120000  # was 0
$ws.Range("K75").Value2 = 0  # was 50000
$ws.Range("L75").Value2 = 120000  # was 0
$ws.Range("M75").ClearContents() | Out-Null  # was -49064
$ws.Range("N75").Value2 = -121872  # new cell
$ws.Range("H76").Value2 = 0  # was 69999.5
$ws.Range("J76").Value2 = 0  # was 69999.5
$ws.Range("L76").Value2 = 0  # was 69999.5
$ws.Range("N76").ClearContents() | Out-Null  # was -70629.5
$ws.Range("H78").Value2 = 120000  # was 50000
$ws.Range("I78").Value2 = 0  # was 50000
$ws.Range("J78").Value2 = 120000  # was 0
$ws.Range("K78").Value2 = 0  # was 150000
$ws.Range("L78").Value2 = 360000  # was 0
$ws.Range("M78").ClearContents() | Out-Null  # was -145320
$ws.Range("N78").Value2 = -369360  # new cell
$ws.Range("H79").Value2 = 0  # was 69999.5
$ws.Range("J79").Value2 = 0  # was 69999.5
$ws.Range("L79").Value2 = 0  # was 69999.5
$ws.Range("N79").ClearContents() | Out-Null  # was -72183.5
$ws.Range("H132").Value2 = 11341942  # was 12081557
$ws.Range("I132").Value2 = 1738150.5  # was 1917836.9
$ws.Range("K132").Value2 = 5214451.5  # was 5753510.699999999
$ws.Range("M132").Value2 = -5211921.5  # was -5750980.699999999
$ws.Range("H136").Value2 = 7892.5728  # was 8115.4194
$ws.Range("I136").Value2 = 3622.6538  # was 3966.7827
$ws.Range("K136").Value2 = 10867.9614  # was 11900.3481
$ws.Range("M136").Value2 = -8317.9614  # was -9350.348100000001
